$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Determine the last used row on this sheet (header row always present).
    $lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11

    # Header: H used to be the "has stage" boolean column, I was "Link", J was "Cover".
    # New layout: H = Link, I = Cover, J removed.
    $ws.Range("H1").Value = "Link"
    $ws.Range("I1").Value = "Cover"

    if ($lastRow -ge 2) {
        for ($r = 2; $r -le $lastRow; $r++) {
            $linkVal = $ws.Cells.Item($r, 9).Value()   # old column I (Link)
            $coverVal = $ws.Cells.Item($r, 10).Value()  # old column J (Cover)
            $ws.Cells.Item($r, 8).Value = $linkVal      # new column H = Link
            $ws.Cells.Item($r, 9).Value = $coverVal     # new column I = Cover
        }
    }

    # Remove the now-duplicate Cover column (old J).
    $ws.Columns.Item(10).Delete()
}

# Refine the "want to go" counts (column F) on sheets that carry event rows.
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Range("F2").Value = 301
$sheet1.Range("F3").Value = 234
$sheet1.Range("F4").Value = 43
$sheet1.Range("F5").Value = 270

$sheet4 = $wb.Worksheets.Item(4)
$sheet4.Range("F2").Value = 301
$sheet4.Range("F3").Value = 234
$sheet4.Range("F4").Value = 43
$sheet4.Range("F5").Value = 270
